# fix: renamed sheet names for working with test_regression_core.py
#
# Mirrors the authoring diff:
#  - Sheet1 "Уравнение регрессии" -> "Лист1"
#  - Sheet2 "Проверка через функции Excel" -> "Расчет"
#  - Sheet1 view zoom 75% -> 100% (normal-zoom persisted)
#  - Sheet2 scroll/selection moved from A48/O73 to A82/L100
#  - Sheet2 column A width set to fit its (now wider) content
#  - Sheet2 O2:O27 / Q2:Q27 re-entered as one block each, which Excel
#    stores as shared formulas (t="shared") instead of per-cell formulas

$wb = $excel.ActiveWorkbook

$wsRegression = $wb.Worksheets.Item(1)
$wsCheck      = $wb.Worksheets.Item(2)

# --- rename sheets -------------------------------------------------------
$wsRegression.Name = "Лист1"
$wsCheck.Name      = "Расчет"

# --- sheet1: zoom 75 -> 100 -----------------------------------------------
$wsRegression.Activate()
$wsRegression.Range("O21").Select()
$excel.ActiveWindow.Zoom = 100

# --- sheet2: re-enter the O/Q columns as one block each (-> shared formulas)
$wsCheck.Range("O2:O27").Formula = "=L2-N2"
$wsCheck.Range("Q2:Q27").Formula = "=ABS(O2)/L2"

# --- sheet2: widen column A to fit its content ----------------------------
$wsCheck.Columns.Item(1).ColumnWidth = 27

# --- sheet2: scroll the viewport and move the selection -------------------
$wsCheck.Activate()
$excel.ActiveWindow.ScrollRow = 82
$excel.ActiveWindow.ScrollColumn = 1
$wsCheck.Range("L100").Select()
